$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.594.13'
$ws.Range('E2').Value = '  -0.14%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.730.62'
$ws.Range('E3').Value = '  -0.97%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9990'
$ws.Range('E4').Value = '  -0.06%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '246.37'
$ws.Range('E5').Value = '  -0.38%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9993'
$ws.Range('E6').Value = '  -0.09%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4818'
$ws.Range('E7').Value = '  +0.23%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2676'
$ws.Range('E8').Value = '  -1.04%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06185'
$ws.Range('E9').Value = '  -1.21%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.730.28'
$ws.Range('E10').Value = '  -0.89%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07165'
$ws.Range('E11').Value = '  +0.79%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.64'
$ws.Range('E12').Value = '  -1.22%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.6115'
$ws.Range('E13').Value = '  -1.04%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.541'
$ws.Range('E14').Value = '  +0.69%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '77.31'
$ws.Range('E15').Value = '  -0.01%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.000'
$ws.Range('E16').Value = '  +0.01%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.586.89'
$ws.Range('E17').Value = '  -0.16%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9993'
$ws.Range('E18').Value = '  -0.10%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000006969'
$ws.Range('E19').Value = '  +0.87%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.56'
$ws.Range('E20').Value = '  -1.35%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.953.67'
$ws.Range('E21').Value = '  -0.87%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.518'
$ws.Range('E22').Value = '  -2.98%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.813'
$ws.Range('E23').Value = '  -0.57%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.252'
$ws.Range('E24').Value = '  -1.96%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '137.31'
$ws.Range('E25').Value = '  +0.74%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.36'
$ws.Range('E26').Value = '  -0.86%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.781'
$ws.Range('E27').Value = '  -2.50%  '

$ws.Range('E28').Value = '  -0.85%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '108.46'
$ws.Range('E29').Value = '  +0.81%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.974'
$ws.Range('E30').Value = '  -1.35%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08014'
$ws.Range('E31').Value = '  +1.43%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.692'
$ws.Range('E32').Value = '  -2.11%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04527'
$ws.Range('E33').Value = '  -1.15%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.611'
$ws.Range('E34').Value = '  -0.13%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.006'
$ws.Range('E35').Value = '  +0.50%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6346'
$ws.Range('E36').Value = '  -0.16%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9003'
$ws.Range('E37').Value = '  -5.98%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.050'
$ws.Range('E38').Value = '  +3.62%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.367'
$ws.Range('E39').Value = '  -4.47%  '

$ws.Range('E40').Value = '  -0.02%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '103.13'
$ws.Range('E41').Value = '  -10.12%  '

$ws.Range('E42').Value = '  -1.13%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.463'
$ws.Range('E43').Value = '  -4.36%  '

$ws.Range('B44').Value = 'Aptos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '7.169'
$ws.Range('E44').Value = '  +5.97%  '

$ws.Range('B45').Value = 'TheSandbox'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.3898'
$ws.Range('E45').Value = '  -0.60%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1186'
$ws.Range('E46').Value = '  -1.67%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05385'
$ws.Range('E47').Value = '  +1.11%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.879'
$ws.Range('E48').Value = '  -1.08%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '30.70'
$ws.Range('E49').Value = '  -0.53%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.253'
$ws.Range('E50').Value = '  +1.55%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '51.43'
$ws.Range('E51').Value = '  -0.66%  '
